$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ------------------------------------------------------------------
# 1) Refresh the "time_taken" timestamps on the existing "data" sheet
# ------------------------------------------------------------------
$dataSheet.Range("F2").Value = "2021-10-05 14:22:57.222803"
$dataSheet.Range("F3").Value = "2021-10-05 14:22:57.222815"
$dataSheet.Range("F4").Value = "2021-10-05 14:22:57.222819"
$dataSheet.Range("F5").Value = "2021-10-05 14:22:57.222822"
$dataSheet.Range("F6").Value = "2021-10-05 14:22:57.222825"

# ------------------------------------------------------------------
# 2) Add a new "metadata" worksheet right after "data"
# ------------------------------------------------------------------
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "metadata"

# Helper cells used to clone the bold/bordered header style (s="1")
# and the matching index-column style already used on the "data" sheet.
$headerStyleSrc = $dataSheet.Range("B1")
$indexStyleSrc = $dataSheet.Range("A2")

# ---- Header row (row 1) ----
$headers = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
$cols = @("B", "C", "D", "E", "F", "G")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $meta.Range($cols[$i] + "1")
    $headerStyleSrc.Copy()
    $cell.PasteSpecial(-4122)
    $cell.Value = $headers[$i]
}

# ---- Data row (row 2) ----
$a2 = $meta.Range("A2")
$indexStyleSrc.Copy()
$a2.PasteSpecial(-4122)
$a2.Value = 0

$meta.Range("B2").Value = "Thrombocythaemia"
$meta.Range("C2").Value = 945

# data_version "1.2" must be stored as text, not a number, so force a
# text number-format before writing it, then restore the default style
# (mirrors how the author's original cell has no explicit style).
$d2 = $meta.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "1.2"
$d2.Style = "Normal"

$meta.Range("E2").Value = "2020-09-30T10:16:19.920144Z"
$meta.Range("F2").Value = "2021-10-05 14:22:57.218149"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/945/?format=json"
